# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.045.19"
$ws.Range("E2").Value = "  -4.92%  "

$ws.Range("D3").Value = "2.223.48"
$ws.Range("E3").Value = "  -6.01%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'319.32"
$ws.Range("E5").Value = "  -3.28%  "

$ws.Range("D6").Value = "'98.50"
$ws.Range("E6").Value = "  -8.75%  "

$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  -8.92%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  -8.45%  "

$ws.Range("D10").Value = "'36.56"
$ws.Range("E10").Value = "  -11.11%  "

$ws.Range("D11").Value = "'54.01"
$ws.Range("E11").Value = "  -3.80%  "

$ws.Range("D12").Value = "'0.0826"
$ws.Range("E12").Value = "  -10.09%  "

$ws.Range("D13").Value = "'7.62"
$ws.Range("E13").Value = "  -10.03%  "

$ws.Range("E14").Value = "  -2.57%  "

$ws.Range("D15").Value = "2.559.51"
$ws.Range("E15").Value = "  -6.43%  "

$ws.Range("E16").Value = "  -12.72%  "

$ws.Range("D17").Value = "'14.27"
$ws.Range("E17").Value = "  -7.31%  "

$ws.Range("D18").Value = "2.219.11"
$ws.Range("E18").Value = "  -6.23%  "

$ws.Range("D19").Value = "42.923.56"
$ws.Range("E19").Value = "  -5.23%  "

$ws.Range("D20").Value = "'13.71"
$ws.Range("E20").Value = "  -9.95%  "

$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -9.65%  "

$ws.Range("D22").Value = "'6.52"
$ws.Range("E22").Value = "  -10.65%  "

$ws.Range("E23").Value = "  -12.15%  "

$ws.Range("D24").Value = "'65.15"
$ws.Range("E24").Value = "  -10.91%  "

$ws.Range("D25").Value = "'235.93"
$ws.Range("E25").Value = "  -9.44%  "

$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  -5.00%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "  -12.00%  "

$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  -3.87%  "

$ws.Range("E31").Value = "  -14.32%  "

$ws.Range("D32").Value = "'36.32"
$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("D33").Value = "'20.24"
$ws.Range("E33").Value = "  -9.44%  "

$ws.Range("D34").Value = "'0.0861"
$ws.Range("E34").Value = "  -10.25%  "

$ws.Range("D35").Value = "'155.19"
$ws.Range("E35").Value = "  -7.50%  "

$ws.Range("E36").Value = "  -6.41%  "

$ws.Range("E37").Value = "  -1.18%  "

$ws.Range("E38").Value = "  -8.29%  "

$ws.Range("E39").Value = "  -5.83%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'4.36"
$ws.Range("E40").Value = "  -7.93%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.105"
$ws.Range("E41").Value = "  -10.59%  "

$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = "  -8.86%  "

$ws.Range("E43").Value = "  -10.16%  "

$ws.Range("D44").Value = "'13.97"
$ws.Range("E44").Value = "  +7.83%  "

$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("D46").Value = "1.734.22"
$ws.Range("E46").Value = "  -7.94%  "

$ws.Range("D47").Value = "'0.202"
$ws.Range("E47").Value = "  -12.34%  "

$ws.Range("D48").Value = "'83.96"
$ws.Range("E48").Value = "  -13.75%  "

$ws.Range("D49").Value = "'8.89"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("D50").Value = "'5.25"
$ws.Range("E50").Value = "  -13.92%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'73.07"
$ws.Range("E51").Value = "  -15.10%  "

